$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.649.30"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.778.88"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("E7").Value = "  -3.25%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("E11").Value = "  +4.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "3.217.42"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").Value = "2.777.45"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.928"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "51.665.06"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.166"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0445"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0837"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.27%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  -5.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.99%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.15%  "
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.39%  "
$ws.Range("D46").Value = "2.083.44"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.942"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.42%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.32%  "
